# PlayerInfo.xlsx - "Move Animation and Card Animation Change" edit
#
# Data changes (Sheet1):
#   - T4  (Zhouzhou / HitDecreaseRate) : 0.9 -> 1
#   - R5  (Timbuktu / NormalAttackCard): "100|101|102"         -> "200|201|202"
#   - Q5  (Timbuktu / InitialCards)    : "103|104|105|106|107" -> "203|204|205|206|207"
#   - Active selection moves to B6
#   - Workbook window position (xWindow/yWindow) moves to 1170,1170
#
# (the new Q5/R5 strings land as two brand-new shared-string entries,
#  exactly as in the target diff: uniqueCount 57 -> 59)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- data edits -------------------------------------------------------
$ws.Range("T4").Value = 1
$ws.Range("R5").Value = "200|201|202"
$ws.Range("Q5").Value = "203|204|205|206|207"

# -- selection ----------------------------------------------------------
$ws.Range("B6").Select() | Out-Null

# -- window position (cosmetic, matches the diff's workbookView change) -
$win = $wb.Windows.Item(1)
$win.Left = 1170
$win.Top = 1170
